# Generate Report for Handoff
# Replace the old UUID-based file name references with the new ones and
# refresh the handoff/generate timestamps, mirroring a fresh localization
# report generation.

$wb = $excel.ActiveWorkbook

$oldGuid = "01e23798-305e-4e53-88de-5936596e6622"
$newGuid = "84b984e4-3f7c-4cbb-a06c-fd148656469f"

$oldZhXlf = "01e23798-305e-4e53-88de-5936596e6622.59788987651c5f609495416330d5c15690e839b6.zh-cn.xlf"
$newZhXlf = "84b984e4-3f7c-4cbb-a06c-fd148656469f.4d7d596b9d962effe887162d476bbe704ceca03d.zh-cn.xlf"

$oldDeXlf = "01e23798-305e-4e53-88de-5936596e6622.59788987651c5f609495416330d5c15690e839b6.de-de.xlf"
$newDeXlf = "84b984e4-3f7c-4cbb-a06c-fd148656469f.4d7d596b9d962effe887162d476bbe704ceca03d.de-de.xlf"

$newGenerateDate = "2016-08-30 09:23:37"
$newZhHandoffDate = "2016-08-30 09:23:32"

$urlPrefix = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a37150697238737ca70a2c4a9607767f53df178c/e2e/"

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A2").Value = "$newGuid.md"

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add(
    $ws.Range("B2"),
    "$urlPrefix$newGuid.md",
    [System.Type]::Missing,
    [System.Type]::Missing,
    "e2e\$newGuid.md"
) | Out-Null

$ws.Range("G2").Value = $newGenerateDate
$ws.Range("G2").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add(
    $ws.Range("A2"),
    "$urlPrefix$newGuid.md",
    [System.Type]::Missing,
    [System.Type]::Missing,
    "$newGuid.md"
) | Out-Null

$ws.Range("G2").Value = $newZhXlf
$ws.Range("H2").Value = $newZhHandoffDate
$ws.Range("H2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("K2").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add(
    $ws.Range("A2"),
    "$urlPrefix$newGuid.md",
    [System.Type]::Missing,
    [System.Type]::Missing,
    "$newGuid.md"
) | Out-Null

$ws.Range("G2").Value = $newDeXlf
$ws.Range("H2").Value = $newGenerateDate
$ws.Range("H2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("K2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
